# Apply the cryptos list update for this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.785.66'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.090.35'
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.74'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.53'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  -4.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0768'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.109'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.882'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.12%  '
$ws.Range("D15").Value = '2.394.08'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("E16").Value = '  -3.95%  '
$ws.Range("D17").Value = '2.084.90'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").Value = '36.778.70'
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("E22").Value = '  +1.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.56%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '21.02'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.54%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.72%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.124'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0610'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.90%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0839'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.74%  '
$ws.Range("E39").Value = '  -3.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.52%  '
$ws.Range("E41").Value = '  +1.50%  '
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0955'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.68%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.339.71'
$ws.Range("E47").Value = '  +4.40%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.93%  '
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("D51").Value = '2.278.38'
$ws.Range("E51").Value = '  +1.48%  '
